$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N3").Value = 17
$ws.Range("Y3").Value = 9
$ws.Range("Z3").Value = 8
$ws.Range("AB3").Value = 26
$ws.Range("AE3").Value = 21
$ws.Range("AL3").Value = 67
$ws.Range("AQ3").Value = 15
$ws.Range("AV3").Value = 51
$ws.Range("AX3").Value = 10
$ws.Range("Q4").Value = 1.73
$ws.Range("R4").Value = 2.08
$ws.Range("Q5").Value = 2.38
$ws.Range("R5").Value = 1.57
$ws.Range("W6").Value = 8
$ws.Range("AA6").Value = 15
$ws.Range("AL6").Value = 29
$ws.Range("BA6").Value = 67
$ws.Range("G7").Value = 3.9
$ws.Range("H7").Value = 3.8
$ws.Range("I7").Value = 1.83
$ws.Range("K7").Value = 2.3
$ws.Range("L7").Value = 2.4
$ws.Range("S7").Value = 1.33
$ws.Range("T7").Value = 3.25
$ws.Range("AD7").Value = 7.5
$ws.Range("AL7").Value = 13
$ws.Range("AM7").Value = 21
$ws.Range("AT7").Value = 3.25
$ws.Range("AY7").Value = 9.5
$ws.Range("BA7").Value = 29
$ws.Range("K8").Value = 2.2
$ws.Range("L8").Value = 3
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 11
$ws.Range("O8").Value = 1.25
$ws.Range("P8").Value = 4
$ws.Range("Q8").Value = 1.83
$ws.Range("R8").Value = 2.03
$ws.Range("U8").Value = 1.67
$ws.Range("V8").Value = 2.1
$ws.Range("X8").Value = 15
$ws.Range("AA8").Value = 23
$ws.Range("AC8").Value = 12
$ws.Range("AH8").Value = 9
$ws.Range("AJ8").Value = 9.5
$ws.Range("AO8").Value = 17
$ws.Range("AY8").Value = 13
$ws.Range("M9").Value = 1.05
$ws.Range("N9").Value = 11
$ws.Range("Q9").Value = 1.98
$ws.Range("R9").Value = 1.88
$ws.Range("AW9").Value = 126
$ws.Range("G12").Value = 1.62
$ws.Range("H12").Value = 3.65
$ws.Range("I12").Value = 5
$ws.Range("O12").Value = 1.26
$ws.Range("P12").Value = 3.15
$ws.Range("S12").Value = 1.39
$ws.Range("U12").Value = 1.78
$ws.Range("V12").Value = 1.82
$ws.Range("W12").Value = 6.8
$ws.Range("X12").Value = 7.6
$ws.Range("Y12").Value = 8
$ws.Range("Z12").Value = 12
$ws.Range("AB12").Value = 26
$ws.Range("AC12").Value = 10.25
$ws.Range("AD12").Value = 7.2
$ws.Range("AE12").Value = 16
$ws.Range("AH12").Value = 13.5
$ws.Range("AI12").Value = 30
$ws.Range("AN12").Value = 3.4
$ws.Range("AO12").Value = 7.9
$ws.Range("AP12").Value = 17.5
$ws.Range("AR12").Value = 60
$ws.Range("AS12").Value = 250
$ws.Range("AU12").Value = 7.6
$ws.Range("AV12").Value = 75
